$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (A3:A21) holds plain text values such as "28/07/2022".
# This change swaps the separator from "/" to "-" for every date in that
# column. Several of the resulting strings (e.g. "01-08-2022") look like
# valid dd-mm-yyyy dates, and Excel's automatic type inference would
# otherwise silently convert them into date serial numbers instead of
# leaving them as plain text. To avoid that, the whole date range is
# temporarily forced to Text format while the new values are written, then
# its style is restored to Normal (the workbook's original/default style
# for these cells) so no visible formatting change is left behind.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

$ws.Range("A3").Value  = "28-07-2022"
$ws.Range("A4").Value  = "01-08-2022"
$ws.Range("A5").Value  = "04-08-2022"
$ws.Range("A6").Value  = "08-08-2022"
$ws.Range("A7").Value  = "11-08-2022"
$ws.Range("A8").Value  = "15-08-2022"
$ws.Range("A9").Value  = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

$dateRange.Style = "Normal"

# Attendance-count corrections for row 3 (28-07-2022) and row 10 (22-08-2022).
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0
